$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new year columns (L = 2022, M = 2023) -------------------------
# First copy the formatting of the existing J:K block (rows 3-10) into the
# new L:M columns so the new cells inherit the same styles (borders, fonts,
# number formats) as the rest of the table.
$ws.Range("J3:K10").Copy() | Out-Null
$ws.Range("L3:M10").PasteSpecial(-4122) | Out-Null

# Year header row
$ws.Range("L4").Value = 2022
$ws.Range("M4").Value = 2023

# Data rows
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 700

$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = 6

$ws.Range("L7").Value = 23

# M7 uses the same "dash" style as the other placeholder cells (e.g. K8)
# rather than the plain style copied from the rest of row 7, so copy that
# cell's formatting across before setting its value.
$ws.Range("K8").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Value = "-"

$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 5

$ws.Range("L9").Value = 23
$ws.Range("M9").Value = 21

$ws.Range("L10").Value = 172
$ws.Range("M10").Value = 143

# --- Footnote row: make A11:C11 use a smaller (8pt) Times New Roman font ---
$footnote = $ws.Range("A11:C11")
$footnote.Font.Name = "Times New Roman"
$footnote.Font.Size = 8
$footnote.VerticalAlignment = -4108

# --- Clear the stale selection and reset it to A1 --------------------------
$ws.Range("A1").Select() | Out-Null

# --- Page setup (printer settings) -----------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 256
$ps.Orientation = 1
